$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# New ballot row (row 20): Jon Heyman, source Fancred, dated 2018-12-06
$ws.Range("A20").Value = "Jon Heyman"
$ws.Range("B20").Value = "x"
$ws.Range("C20").Value = "x"
$ws.Range("E20").Value = "x"
$ws.Range("G20").Value = "x"
$ws.Range("I20").Value = "x"
$ws.Range("J20").Value = "x"
$ws.Range("K20").Value = "x"
$ws.Range("O20").Value = "x"
$ws.Range("P20").Value = "x"
$ws.Range("Q20").Value = "x"
$ws.Range("AK20").Value = 10
$ws.Range("AL20").Value = "Fancred"

# Copy the date format from the row above so we reuse the existing
# built-in date style instead of minting a new custom number format.
$ws.Range("AM19").Copy()
$ws.Range("AM20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("AM20").Value = 43440

# Update the active selection to match the new cursor position.
$ws.Range("C21").Select() | Out-Null
